$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 14.11239547175637
$ws.Range("F3").Value = 13.09487473480318
$ws.Range("F4").Value = 8.22111200880744
$ws.Range("F5").Value = 7.429121582096163
$ws.Range("F6").Value = 6.324528075904071
$ws.Range("F7").Value = 6.243826188088984
$ws.Range("F10").Value = 5.141087836715284
$ws.Range("F11").Value = 3.417079858592328
$ws.Range("F12").Value = 1.153463192899035
$ws.Range("F13").Value = 0.3269558257719956

# Rows 8 and 9 swap B, C, D, F, G values
$ws.Range("B8").Value = 32
$ws.Range("C8").Value = "60bf9943e4e04642d4634ecc"
$ws.Range("D8").Value = "Jamarii"
$ws.Range("F8").Value = 5.27722767756892
$ws.Range("G8").Value = "Black or African American"

$ws.Range("B9").Value = 33
$ws.Range("C9").Value = "60b322994d0b901954690036"
$ws.Range("D9").Value = "Brennan"
$ws.Range("F9").Value = 5.186042016282854
$ws.Range("G9").Value = "White"
